$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the funny header column names (row 1)
$ws.Range("A1").Value = "section1"
$ws.Range("B1").Value = "name__"
$ws.Range("C1").Value = "des"
$ws.Range("D1").Value = "data_typess"
$ws.Range("E1").Value = "necessary"
$ws.Range("F1").Value = "spec2"

# Move the active selection to A2 as the active cell
$ws.Range("A2").Select()
